# MVP for buy n get m at x%
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: add a "Limit" value (column H) for the existing "nforx" specialty ---
$ws.Range("H11").Value = 6

# --- Row 17 (Banana): turn on a new "buy n get m at x%" specialty ---
$ws.Range("E17").Value = 0.49          # Markdown
$ws.Range("F17").Value = $true         # Has Specialty
$ws.Range("G17").Value = "nmatx"       # Type of Specialty (new shared string)
$ws.Range("H17").Value = 6             # Limit
$ws.Range("I17").Value = 2             # Specialty Variable 1 (buy n)
$ws.Range("J17").Value = 1             # Specialty Variable 2 (get m)
$ws.Range("K17").Value = 0.5           # Specialty Variable 3 (at x%)

# --- Window / view bookkeeping to mirror the author's on-screen state ---
try { $excel.ActiveWindow.Left = 11620 } catch {}
try { $excel.ActiveWindow.Top = 1900 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 3 } catch {}
try { $excel.ActiveWindow.ScrollRow = 1 } catch {}

$ws.Range("H16").Select()

Write-Output "applied buy-n-get-m-at-x% MVP changes"
